# Auto-generated edit script applying the Ifrit_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 76926000
$ws.Range("I137").Value = 1974.75
$ws.Range("J137").Value = 111114456
$ws.Range("K137").Value = 5924.25
$ws.Range("L137").Value = 333343368
$ws.Range("M137").Value = -3374.25
$ws.Range("N137").Value = -333348468

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 4132.5
$ws.Range("I21").Value = 959
$ws.Range("J21").Value = 20000
$ws.Range("K21").Value = 959
$ws.Range("L21").Value = 20000
$ws.Range("M21").Value = -585
$ws.Range("N21").Value = -20748

$ws.Range("H74").Value = 4085.225
$ws.Range("I74").Value = 1056.3334
$ws.Range("J74").Value = 5383.3213
$ws.Range("K74").Value = 1056.3334
$ws.Range("L74").Value = 5383.3213
$ws.Range("M74").Value = -182.3334
$ws.Range("N74").Value = -7131.3213

$ws.Range("H77").Value = 4085.225
$ws.Range("I77").Value = 1056.3334
$ws.Range("J77").Value = 5383.3213
$ws.Range("K77").Value = 5281.666999999999
$ws.Range("L77").Value = 26916.6065
$ws.Range("M77").Value = -913.6669999999995
$ws.Range("N77").Value = -35652.60649999999

$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws.Range("H109").Value = 40188.5
$ws.Range("J109").Value = 40188.5
$ws.Range("L109").Value = 40188.5
$ws.Range("N109").Value = -42962.5

$ws.Range("H132").Value = 3491.4736
$ws.Range("I132").Value = 3317.1428
$ws.Range("J132").Value = 3979.6
$ws.Range("K132").Value = 9951.428400000001
$ws.Range("L132").Value = 11938.8
$ws.Range("M132").Value = -7421.428400000001
$ws.Range("N132").Value = -16998.8

$ws.Range("H133").Value = 122000
$ws.Range("J133").Value = 122000
$ws.Range("L133").Value = 122000
$ws.Range("N133").Value = -127060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 663.2222
$ws.Range("I80").Value = 163.125
$ws.Range("J80").Value = 1063.3
$ws.Range("K80").Value = 163.125
$ws.Range("L80").Value = 1063.3
$ws.Range("M80").Value = 834.875
$ws.Range("N80").Value = -3059.3

$ws.Range("H83").Value = 663.2222
$ws.Range("I83").Value = 163.125
$ws.Range("J83").Value = 1063.3
$ws.Range("K83").Value = 815.625
$ws.Range("L83").Value = 5316.5
$ws.Range("M83").Value = 4176.375
$ws.Range("N83").Value = -15300.5

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 27309.363
$ws.Range("J92").Value = 27309.363
$ws.Range("L92").Value = 27309.363
$ws.Range("N92").Value = -32301.363

$ws.Range("H94").Value = 1408
$ws.Range("I94").Value = 1612
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 1612
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -1161
$ws.Range("N94").Value = -1902

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1773.75
$ws.Range("I9").Value = 200
$ws.Range("J9").Value = 1998.5714
$ws.Range("K9").Value = 600
$ws.Range("L9").Value = 5995.7142
$ws.Range("M9").Value = -376
$ws.Range("N9").Value = -6443.7142

$ws.Range("H15").Value = 635
$ws.Range("J15").Value = 833.3333
$ws.Range("L15").Value = 2499.9999
$ws.Range("N15").Value = -2779.9999

$ws.Range("H16").Value = 30001
$ws.Range("I16").Value = 30001
$ws.Range("K16").Value = 90003
$ws.Range("M16").Value = -89830

$ws.Range("H19").Value = 2187
$ws.Range("J19").Value = 2187
$ws.Range("L19").Value = 6561
$ws.Range("N19").Value = -6909

$ws.Range("H113").Value = 496.30234
$ws.Range("I113").Value = 473.33334
$ws.Range("J113").Value = 525.3158
$ws.Range("K113").Value = 1420.00002
$ws.Range("L113").Value = 1575.9474
$ws.Range("M113").Value = 749.9999800000001
$ws.Range("N113").Value = -5915.9474

$ws.Range("H122").Value = 10753890
$ws.Range("I122").Value = 16667124
$ws.Range("J122").Value = 2555.3635
$ws.Range("K122").Value = 150004116
$ws.Range("L122").Value = 22998.2715
$ws.Range("M122").Value = -150001666
$ws.Range("N122").Value = -27898.2715

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 525.6667
$ws.Range("I107").Value = 325.92307
$ws.Range("J107").Value = 761.7273
$ws.Range("K107").Value = 325.92307
$ws.Range("L107").Value = 761.7273
$ws.Range("M107").Value = 1594.07693
$ws.Range("N107").Value = -4601.7273

$ws.Range("H123").Value = 17326
$ws.Range("J123").Value = 17326
$ws.Range("L123").Value = 17326
$ws.Range("N123").Value = -22226

$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -50120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 953
$ws.Range("I24").Value = 953
$ws.Range("K24").Value = 953
$ws.Range("M24").Value = -610

$ws.Range("H61").Value = 626.1818
$ws.Range("I61").Value = 570.8889
$ws.Range("K61").Value = 570.8889
$ws.Range("M61").Value = -368.8889

$ws.Range("H113").Value = 626.1818
$ws.Range("I113").Value = 570.8889
$ws.Range("K113").Value = 570.8889
$ws.Range("M113").Value = 1599.1111

$ws.Range("H132").Value = 6382.185
$ws.Range("I132").Value = 9461.066000000001
$ws.Range("J132").Value = 2533.5833
$ws.Range("K132").Value = 28383.198
$ws.Range("L132").Value = 7600.749899999999
$ws.Range("M132").Value = -25853.198
$ws.Range("N132").Value = -12660.7499

$ws.Range("H133").Value = 59999
$ws.Range("J133").Value = 59999
$ws.Range("L133").Value = 59999
$ws.Range("N133").Value = -65059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 50000
$ws.Range("J29").Value = 50000
$ws.Range("L29").Value = 50000
$ws.Range("N29").Value = -50580

$ws.Range("H123").Value = 44320
$ws.Range("J123").Value = 44320
$ws.Range("L123").Value = 44320
$ws.Range("N123").Value = -54120

$ws.Range("H132").Value = 1520.6316
$ws.Range("I132").Value = 773.625
$ws.Range("J132").Value = 2063.9092
$ws.Range("K132").Value = 2320.875
$ws.Range("L132").Value = 6191.7276
$ws.Range("M132").Value = 209.125
$ws.Range("N132").Value = -11251.7276
